$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

# The "Employment" section's last row (row 60, "Employment impact") is being
# removed entirely. Deleting the whole row shifts everything below up by one
# and shrinks/re-numbers the merged cell ranges accordingly, matching the
# target diff (dimension A1:N152 -> A1:N151).
$ws.Rows(60).Delete()
